# Applies the Serpinc1-Sdc2 LR-pairs update (Natmi following Dr Hou advice):
# recompute rows 2-9 with updated values and append new rows 10-17.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = @("ECs", "Serpinc1", "Sdc2", "ECs", 3, 1, 4.439248, 13.317744, 0.3193177756555054, 0.3193177756555054, 2, 0.6666666666666666, 0.6646083333333334, 1.993825, 0.005515555660921567, 0.005515555660921567, 2.950361214533333, 26.5532509308, 0.001761214965149606, 0.001761214965149606)
$row3 = @("ECs", "Serpinc1", "Sdc2", "FAPs", 3, 1, 4.439248, 13.317744, 0.3193177756555054, 0.3193177756555054, 3, 1, 85.826024, 257.478072, 0.7122664414292983, 0.7122664414292983, 381.003005389952, 3429.027048509568, 0.2274393357512658, 0.2274393357512658)
$row4 = @("ECs", "Serpinc1", "Sdc2", "M2", 3, 1, 4.439248, 13.317744, 0.3193177756555054, 0.3193177756555054, 1, 0.3333333333333333, 0.009315666666666667, 0.027947, 0.00007731031261809587, 0.00007731031261809588, 0.04135455461866667, 0.372190991568, 0.00002468655706044212, 0.00002468655706044213)
$row5 = @("ECs", "Serpinc1", "Sdc2", "sCs", 3, 1, 4.439248, 13.317744, 0.3193177756555054, 0.3193177756555054, 3, 1, 33.99712866666667, 101.991386, 0.2821406925971621, 0.2821406925971621, 150.9216854392427, 1358.295168953184, 0.0900925383820295, 0.0900925383820295)
$row6 = @("FAPs", "Serpinc1", "Sdc2", "ECs", 3, 1, 3.893859666666666, 11.681579, 0.2800876651799331, 0.2800876651799331, 2, 0.6666666666666666, 0.6646083333333334, 1.993825, 0.005515555660921567, 0.005515555660921567, 2.587891583297222, 23.291024249675, 0.001544839107237484, 0.001544839107237484)
$row7 = @("FAPs", "Serpinc1", "Sdc2", "FAPs", 3, 1, 3.893859666666666, 11.681579, 0.2800876651799331, 0.2800876651799331, 3, 1, 85.826024, 257.478072, 0.7122664414292983, 0.7122664414292983, 334.1944932039653, 3007.750438835688, 0.1994970445659517, 0.1994970445659517)
$row8 = @("FAPs", "Serpinc1", "Sdc2", "M2", 3, 1, 3.893859666666666, 11.681579, 0.2800876651799331, 0.2800876651799331, 1, 0.3333333333333333, 0.009315666666666667, 0.027947, 0.00007731031261809587, 0.00007731031261809588, 0.03627389870144444, 0.326465088313, 0.00002165366495553319, 0.0000216536649555332)
$row9 = @("FAPs", "Serpinc1", "Sdc2", "sCs", 3, 1, 3.893859666666666, 11.681579, 0.2800876651799331, 0.2800876651799331, 3, 1, 33.99712866666667, 101.991386, 0.2821406925971621, 0.2821406925971621, 132.3800480976104, 1191.420432878494, 0.07902412784178835, 0.07902412784178836)
$row10 = @("M2", "Serpinc1", "Sdc2", "ECs", 3, 1, 1.735205333333333, 5.205616, 0.1248143621049263, 0.1248143621049263, 2, 0.6666666666666666, 0.6646083333333334, 1.993825, 0.005515555660921567, 0.005515555660921567, 1.153231924577778, 10.3790873212, 0.0006884205614721404, 0.0006884205614721405)
$row11 = @("M2", "Serpinc1", "Sdc2", "FAPs", 3, 1, 1.735205333333333, 5.205616, 0.1248143621049263, 0.1248143621049263, 3, 1, 85.826024, 257.478072, 0.7122664414292983, 0.7122664414292983, 148.9257745835947, 1340.331971252352, 0.08890108153574369, 0.0889010815357437)
$row12 = @("M2", "Serpinc1", "Sdc2", "M2", 3, 1, 1.735205333333333, 5.205616, 0.1248143621049263, 0.1248143621049263, 1, 0.3333333333333333, 0.009315666666666667, 0.027947, 0.00007731031261809587, 0.00007731031261809588, 0.01616459448355555, 0.145481350352, 0.000009649437353560069, 0.00000964943735356007)
$row13 = @("M2", "Serpinc1", "Sdc2", "sCs", 3, 1, 1.735205333333333, 5.205616, 0.1248143621049263, 0.1248143621049263, 3, 1, 33.99712866666667, 101.991386, 0.2821406925971621, 0.2821406925971621, 58.99199898041956, 530.927990823776, 0.03521521057035688, 0.03521521057035688)
$row14 = @("sCs", "Serpinc1", "Sdc2", "ECs", 3, 1, 3.833976, 11.501928, 0.2757801970596353, 0.2757801970596353, 2, 0.6666666666666666, 0.6646083333333334, 1.993825, 0.005515555660921567, 0.005515555660921567, 2.5480923994, 22.9328315946, 0.001521081027062337, 0.001521081027062337)
$row15 = @("sCs", "Serpinc1", "Sdc2", "FAPs", 3, 1, 3.833976, 11.501928, 0.2757801970596353, 0.2757801970596353, 3, 1, 85.826024, 257.478072, 0.7122664414292983, 0.7122664414292983, 329.054916191424, 2961.494245722816, 0.1964289795763371, 0.1964289795763371)
$row16 = @("sCs", "Serpinc1", "Sdc2", "M2", 3, 1, 3.833976, 11.501928, 0.2757801970596353, 0.2757801970596353, 1, 0.3333333333333333, 0.009315666666666667, 0.027947, 0.00007731031261809587, 0.00007731031261809588, 0.035716042424, 0.321444381816, 0.00002132065324856049, 0.00002132065324856049)
$row17 = @("sCs", "Serpinc1", "Sdc2", "sCs", 3, 1, 3.833976, 11.501928, 0.2757801970596353, 0.2757801970596353, 3, 1, 33.99712866666667, 101.991386, 0.2821406925971621, 0.2821406925971621, 130.344175376912, 1173.097578392208, 0.07780881580298735, 0.07780881580298735)

$allRows = @($row2, $row3, $row4, $row5, $row6, $row7, $row8, $row9, $row10, $row11, $row12, $row13, $row14, $row15, $row16, $row17)

$rowIdx = 2
foreach ($dataRow in $allRows) {
  $colIdx = 1
  foreach ($val in $dataRow) {
    $ws.Cells.Item($rowIdx, $colIdx).Value = $val
    $colIdx = $colIdx + 1
  }
  $rowIdx = $rowIdx + 1
}
